# Commit: redefined prefix "ome" instead of ":" (base prefix name).
#
# The workbook models a Turtle-like prefix table on the "@prefix" sheet and
# uses those prefixes (e.g. ":Image", ":pixels", ...) throughout the other
# sheets. This change renames the base/default prefix from the empty string
# ("") to "ome", so every value that used to start with a bare ":" now
# starts with "ome:".

$wb = $excel.ActiveWorkbook

# 1) @prefix sheet: the base prefix name (was blank) becomes "ome".
$wsPrefix = $wb.Worksheets.Item("@prefix")
$wsPrefix.Range("A1").Value = "ome"

# 2) Image sheet: PropertyURI / datatype columns using the base prefix.
$wsImage = $wb.Worksheets.Item("Image")
$wsImage.Range("E3").Value = "ome:pixels"
$wsImage.Range("F3").Value = "ome:acquisitionDate"
$wsImage.Range("B4").Value = "ome:Image"
$wsImage.Range("E4").Value = "ome:Pixels"

# 3) Pixels sheet.
$wsPixels = $wb.Worksheets.Item("Pixels")
$wsPixels.Range("D3").Value = "ome:pixelType"
$wsPixels.Range("E3").Value = "ome:dimensionOrder"
$wsPixels.Range("F3").Value = "ome:sizeC"
$wsPixels.Range("G3").Value = "ome:sizeT"
$wsPixels.Range("H3").Value = "ome:sizeX"
$wsPixels.Range("I3").Value = "ome:sizeY"
$wsPixels.Range("J3").Value = "ome:sizeZ"
$wsPixels.Range("K3").Value = "ome:channel"
$wsPixels.Range("L3").Value = "ome:metadataOnly"
$wsPixels.Range("B4").Value = "ome:Pixels"
$wsPixels.Range("D4").Value = "ome:PixelType"
$wsPixels.Range("E4").Value = "ome:DimensionOrder"
$wsPixels.Range("K4").Value = "ome:Channel"

# 4) Channel sheet.
$wsChannel = $wb.Worksheets.Item("Channel")
$wsChannel.Range("D3").Value = "ome:color"
$wsChannel.Range("B4").Value = "ome:Channel"
$wsChannel.Range("D4").Value = "ome:Color"

# 5) Color sheet.
$wsColor = $wb.Worksheets.Item("Color")
$wsColor.Range("B4").Value = "ome:Color"
